# Corrects the IFRS financial-summary figures in the "company_list" sheet.
# Rows 2-6 (FY2014-FY2018 actuals) had their D:AJ figures replaced with the
# correct (much smaller, properly-scaled) numbers, and rows 7-9 (the 2019E-2021E
# estimate rows, which held stale/incorrect placeholder figures) are cleared back
# to just their rank/period/label columns (A:C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-6: overwrite columns D through AJ with corrected figures ---

# Row 2
$ws.Range("D2").Value = 4594
$ws.Range("E2").Value = 192
$ws.Range("F2").Value = 192
$ws.Range("G2").Value = 198
$ws.Range("H2").Value = 175
$ws.Range("I2").Value = 172
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 4308
$ws.Range("L2").Value = 2223
$ws.Range("M2").Value = 2085
$ws.Range("N2").Value = 2020
$ws.Range("O2").Value = 65
$ws.Range("P2").Value = 95
$ws.Range("Q2").Value = 222
$ws.Range("R2").Value = -646
$ws.Range("S2").Value = 329
$ws.Range("T2").Value = 612
$ws.Range("U2").Value = -390
$ws.Range("V2").Value = 1106
$ws.Range("W2").Value = 4.18
$ws.Range("X2").Value = 3.81
$ws.Range("Y2").Value = 8.8
$ws.Range("Z2").Value = 4.32
$ws.Range("AA2").Value = 106.62
$ws.Range("AB2").Value = 2007.7
$ws.Range("AC2").Value = 902
$ws.Range("AD2").Value = 7.04
$ws.Range("AE2").Value = 10593
$ws.Range("AF2").Value = 0.6
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 3.15
$ws.Range("AI2").Value = 22.17
$ws.Range("AJ2").Value = 19072280

# Row 3
$ws.Range("D3").Value = 4675
$ws.Range("E3").Value = 169
$ws.Range("F3").Value = 169
$ws.Range("G3").Value = 153
$ws.Range("H3").Value = 120
$ws.Range("I3").Value = 118
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 4594
$ws.Range("L3").Value = 2451
$ws.Range("M3").Value = 2143
$ws.Range("N3").Value = 2076
$ws.Range("O3").Value = 67
$ws.Range("P3").Value = 95
$ws.Range("Q3").Value = 143
$ws.Range("R3").Value = -384
$ws.Range("S3").Value = 154
$ws.Range("T3").Value = 367
$ws.Range("U3").Value = -224
$ws.Range("V3").Value = 1335
$ws.Range("W3").Value = 3.62
$ws.Range("X3").Value = 2.56
$ws.Range("Y3").Value = 5.76
$ws.Range("Z3").Value = 2.69
$ws.Range("AA3").Value = 114.37
$ws.Range("AB3").Value = 2065.92
$ws.Range("AC3").Value = 618
$ws.Range("AD3").Value = 7.6
$ws.Range("AE3").Value = 10884
$ws.Range("AF3").Value = 0.43
$ws.Range("AG3").Value = 120
$ws.Range("AH3").Value = 2.55
$ws.Range("AI3").Value = 19.41
$ws.Range("AJ3").Value = 19072280

# Row 4
$ws.Range("D4").Value = 4951
$ws.Range("E4").Value = 243
$ws.Range("F4").Value = 243
$ws.Range("G4").Value = 228
$ws.Range("H4").Value = 195
$ws.Range("I4").Value = 187
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 4663
$ws.Range("L4").Value = 2374
$ws.Range("M4").Value = 2290
$ws.Range("N4").Value = 2218
$ws.Range("O4").Value = 72
$ws.Range("P4").Value = 95
$ws.Range("Q4").Value = 419
$ws.Range("R4").Value = -335
$ws.Range("S4").Value = -95
$ws.Range("T4").Value = 341
$ws.Range("U4").Value = 78
$ws.Range("V4").Value = 1281
$ws.Range("W4").Value = 4.91
$ws.Range("X4").Value = 3.94
$ws.Range("Y4").Value = 8.73
$ws.Range("Z4").Value = 4.22
$ws.Range("AA4").Value = 103.68
$ws.Range("AB4").Value = 2228.41
$ws.Range("AC4").Value = 983
$ws.Range("AD4").Value = 5.04
$ws.Range("AE4").Value = 11628
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 3.03
$ws.Range("AI4").Value = 15.27
$ws.Range("AJ4").Value = 19072280

# Row 5
$ws.Range("D5").Value = 4668
$ws.Range("E5").Value = 135
$ws.Range("F5").Value = 135
$ws.Range("G5").Value = 133
$ws.Range("H5").Value = 118
$ws.Range("I5").Value = 115
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4689
$ws.Range("L5").Value = 2305
$ws.Range("M5").Value = 2383
$ws.Range("N5").Value = 2291
$ws.Range("O5").Value = 93
$ws.Range("P5").Value = 95
$ws.Range("Q5").Value = 263
$ws.Range("R5").Value = -317
$ws.Range("S5").Value = 122
$ws.Range("T5").Value = 264
$ws.Range("U5").Value = -2
$ws.Range("V5").Value = 1344
$ws.Range("W5").Value = 2.89
$ws.Range("X5").Value = 2.53
$ws.Range("Y5").Value = 5.08
$ws.Range("Z5").Value = 2.53
$ws.Range("AA5").Value = 96.72
$ws.Range("AB5").Value = 2333.27
$ws.Range("AC5").Value = 601
$ws.Range("AD5").Value = 13.86
$ws.Range("AE5").Value = 12011
$ws.Range("AF5").Value = 0.69
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.2
$ws.Range("AI5").Value = 16.64
$ws.Range("AJ5").Value = 19072280

# Row 6
$ws.Range("D6").Value = 4572
$ws.Range("E6").Value = 41
$ws.Range("F6").Value = 41
$ws.Range("G6").Value = -10
$ws.Range("H6").Value = -7
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 4752
$ws.Range("L6").Value = 2441
$ws.Range("M6").Value = 2311
$ws.Range("N6").Value = 2227
$ws.Range("P6").Value = 95
$ws.Range("Q6").Value = 299
$ws.Range("R6").Value = -344
$ws.Range("S6").Value = -6
$ws.Range("T6").Value = 364
$ws.Range("U6").Value = -65
$ws.Range("V6").Value = 1374
$ws.Range("W6").Value = 0.9
$ws.Range("X6").Value = -0.16
$ws.Range("Y6").Value = 0.06
$ws.Range("Z6").Value = -0.15
$ws.Range("AA6").Value = 105.61
$ws.Range("AB6").Value = 2269.04
$ws.Range("AC6").Value = 7
$ws.Range("AD6").Value = 1213.79
$ws.Range("AE6").Value = 11678
$ws.Range("AF6").Value = 0.71
$ws.Range("AG6").Value = 30
$ws.Range("AH6").Value = 0.36
$ws.Range("AI6").Value = 439.78
$ws.Range("AJ6").Value = 19072280

# --- Rows 7-9: clear the stale figures, keep only rank (A), period (B) and label (C) ---
$ws.Range("D7:AI9").ClearContents()

